{"js": "// Replace each three-digit-by-one-digit multiplication prompt in the\n// worksheet table with the new prompt from the commit, matching the\n// exact old text so each cell's formatting (run properties) is kept\n// untouched -- only the run's text content changes.\nconst replacements = [\n  [\"221\u00d78=\", \"959\u00d79=\"],\n  [\"706\u00d78=\", \"639\u00d73=\"],\n  [\"851\u00d79=\", \"758\u00d78=\"],\n  [\"346\u00d76=\", \"336\u00d72=\"],\n  [\"738\u00d75=\", \"353\u00d77=\"],\n  [\"201\u00d76=\", \"558\u00d77=\"],\n  [\"957\u00d76=\", \"943\u00d77=\"],\n  [\"576\u00d73=\", \"311\u00d76=\"],\n  [\"958\u00d72=\", \"324\u00d76=\"],\n  [\"585\u00d79=\", \"319\u00d73=\"],\n  [\"193\u00d78=\", \"137\u00d75=\"],\n  [\"626\u00d76=\", \"457\u00d73=\"],\n  [\"280\u00d75=\", \"741\u00d74=\"],\n  [\"953\u00d79=\", \"802\u00d77=\"],\n  [\"172\u00d72=\", \"906\u00d78=\"],\n  [\"969\u00d72=\", \"606\u00d78=\"],\n  [\"584\u00d72=\", \"933\u00d72=\"],\n  [\"234\u00d75=\", \"782\u00d78=\"],\n  [\"169\u00d77=\", \"223\u00d79=\"],\n  [\"145\u00d78=\", \"228\u00d72=\"],\n  [\"855\u00d74=\", \"169\u00d72=\"],\n  [\"199\u00d78=\", \"508\u00d73=\"],\n  [\"853\u00d77=\", \"350\u00d77=\"],\n  [\"287\u00d74=\", \"471\u00d79=\"],\n  [\"329\u00d76=\", \"105\u00d73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit-by-one-digit multiplication prompt in the\n# worksheet table with the new prompt from the commit. Each old prompt\n# is unique in the document, so a plain Find/Replace-all per pair is\n# safe and leaves the surrounding run formatting (font/size) untouched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"221\u00d78=\", \"959\u00d79=\"),\n    @(\"706\u00d78=\", \"639\u00d73=\"),\n    @(\"851\u00d79=\", \"758\u00d78=\"),\n    @(\"346\u00d76=\", \"336\u00d72=\"),\n    @(\"738\u00d75=\", \"353\u00d77=\"),\n    @(\"201\u00d76=\", \"558\u00d77=\"),\n    @(\"957\u00d76=\", \"943\u00d77=\"),\n    @(\"576\u00d73=\", \"311\u00d76=\"),\n    @(\"958\u00d72=\", \"324\u00d76=\"),\n    @(\"585\u00d79=\", \"319\u00d73=\"),\n    @(\"193\u00d78=\", \"137\u00d75=\"),\n    @(\"626\u00d76=\", \"457\u00d73=\"),\n    @(\"280\u00d75=\", \"741\u00d74=\"),\n    @(\"953\u00d79=\", \"802\u00d77=\"),\n    @(\"172\u00d72=\", \"906\u00d78=\"),\n    @(\"969\u00d72=\", \"606\u00d78=\"),\n    @(\"584\u00d72=\", \"933\u00d72=\"),\n    @(\"234\u00d75=\", \"782\u00d78=\"),\n    @(\"169\u00d77=\", \"223\u00d79=\"),\n    @(\"145\u00d78=\", \"228\u00d72=\"),\n    @(\"855\u00d74=\", \"169\u00d72=\"),\n    @(\"199\u00d78=\", \"508\u00d73=\"),\n    @(\"853\u00d77=\", \"350\u00d77=\"),\n    @(\"287\u00d74=\", \"471\u00d79=\"),\n    @(\"329\u00d76=\", \"105\u00d73=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
